$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 335, shifting existing rows 335:393 down to 339:397
$ws.Rows("335:338").Insert()

# New data block: Clementina prices, date 2021-11-05 (serial 44505), Provincia de Quillota
$newRows = @(
    @{ Row=335; K="Clementina"; L="Especial"; M=350; N=6000; O=6000; P=6000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=600; T=10 },
    @{ Row=336; K="Clementina"; L="Primera";  M=330; N=5000; O=5000; P=5000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=500; T=10 },
    @{ Row=337; K="Clementina"; L="Segunda";  M=300; N=4000; O=4000; P=4000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=400; T=10 },
    @{ Row=338; K="Clementina"; L="Tercera";  M=250; N=3000; O=3000; P=3000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota"; S=300; T=10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44505
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102004
    $ws.Cells.Item($row, 10).Value = "Mandarina"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
